# Add a "الفرع" (Branch) column to the sales report, between the date
# column (D) and the customer-code column (old E, now F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at E; everything from old E..O shifts to F..P.
[void]$ws.Columns("E").Insert()

# 2. Give the new header cell the same look as the other bold/filled header
#    cells (copy formatting from A1), then set its text.
$ws.Range("A1").Copy()
[void]$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "الفرع"

# 3. Size the new column like its neighbour column D.
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# 4. Move the active selection to the new column on row 2.
[void]$ws.Range("E2").Select()

# 5. Re-apply AutoFilter so its range grows from A1:O1 to A1:P1.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:P1").AutoFilter()

# 6. Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=المبيعات!`$A`$1:`$P`$1"
    }
}
